{"js": "// Update the Korean financial-statement header row:\n//   - bold every header label (except \"Year\", which is already bold)\n//   - rewrite \"$1,000\" labels to the more compact \"$K\" style, with a\n//     couple of small wording tweaks to match the updated terminology\nconst replacements = [\n  { from: \"\ub9e4\ucd9c($1,000)\", to: \"\uc218\uc775($K)\" },\n  { from: \"\ud310\ub9e4 \uc81c\ud488 \uc6d0\uac00($1,000)\", to: \"\ud310\ub9e4 \uc81c\ud488 \uc6d0\uac00($K)\" },\n  { from: \"\uc601\uc5c5\ube44\uc6a9($1,000)\", to: \"\uc601\uc5c5 \ube44\uc6a9($K)\" },\n  { from: \"EBITDA($1,000)\", to: \"EBITDA($K)\" },\n  { from: \"\uc774\uc790\ube44\uc6a9($1,000)\", to: \"\uc774\uc790 \ube44\uc6a9($K)\" },\n  { from: \"\uc138\uc804 \uc774\uc775($1,000)\", to: \"\uc138\uc804 \uc774\uc775($K)\" },\n  { from: \"\uc21c\uc218\uc785($1,000)\", to: \"\uc21c\uc774\uc775($K)\" },\n  { from: \"\uc790\uc0b0 \ucd1d\uc561($1,000)\", to: \"\ucd1d \uc790\uc0b0($K)\" },\n  { from: \"\ubd80\ucc44 \ucd1d\uc561($1,000)\", to: \"\ucd1d \ubd80\ucc44($K)\" },\n  { from: \"\uc790\uae30 \uc790\ubcf8($1,000)\", to: \"\uc8fc\uc8fc \uc9c0\ubd84($K)\" },\n];\n\n// Cell whose label text stays the same, but still needs to become bold.\nconst boldOnly = [\"\ucd1d \uc774\uc775\ub960(%)\"];\n\nfor (const { from, to } of replacements) {\n  const results = context.document.body.search(from, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const r of results.items) {\n    r.font.bold = true;\n    r.insertText(to, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n\nfor (const label of boldOnly) {\n  const results = context.document.body.search(label, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const r of results.items) {\n    r.font.bold = true;\n  }\n  await context.sync();\n}\n", "ps1": "# Update the Korean financial-statement header row:\n#   - bold every header label (except \"Year\", which is already bold)\n#   - rewrite \"$1,000\" labels to the more compact \"$K\" style, with a\n#     couple of small wording tweaks to match the updated terminology\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Find = \"\ub9e4\ucd9c(`$1,000)\";         Replace = \"\uc218\uc775(`$K)\" },\n    @{ Find = \"\ud310\ub9e4 \uc81c\ud488 \uc6d0\uac00(`$1,000)\"; Replace = \"\ud310\ub9e4 \uc81c\ud488 \uc6d0\uac00(`$K)\" },\n    @{ Find = \"\uc601\uc5c5\ube44\uc6a9(`$1,000)\";       Replace = \"\uc601\uc5c5 \ube44\uc6a9(`$K)\" },\n    @{ Find = \"EBITDA(`$1,000)\";        Replace = \"EBITDA(`$K)\" },\n    @{ Find = \"\uc774\uc790\ube44\uc6a9(`$1,000)\";       Replace = \"\uc774\uc790 \ube44\uc6a9(`$K)\" },\n    @{ Find = \"\uc138\uc804 \uc774\uc775(`$1,000)\";      Replace = \"\uc138\uc804 \uc774\uc775(`$K)\" },\n    @{ Find = \"\uc21c\uc218\uc785(`$1,000)\";        Replace = \"\uc21c\uc774\uc775(`$K)\" },\n    @{ Find = \"\uc790\uc0b0 \ucd1d\uc561(`$1,000)\";      Replace = \"\ucd1d \uc790\uc0b0(`$K)\" },\n    @{ Find = \"\ubd80\ucc44 \ucd1d\uc561(`$1,000)\";      Replace = \"\ucd1d \ubd80\ucc44(`$K)\" },\n    @{ Find = \"\uc790\uae30 \uc790\ubcf8(`$1,000)\";      Replace = \"\uc8fc\uc8fc \uc9c0\ubd84(`$K)\" }\n)\n\nforeach ($r in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $r.Find\n    $find.Replacement.Text = $r.Replace\n    $find.Replacement.Font.Bold = 1\n    $find.Execute($find.Text, $false, $true, $false, $false, $false, $true, 1, $true, $find.Replacement.Text, 2)\n}\n\n# \"\ucd1d \uc774\uc775\ub960(%)\" keeps its text, but still needs to become bold like the\n# rest of the header row. Replace the text with itself so the Format\n# (bold) change is applied without altering the label.\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"\ucd1d \uc774\uc775\ub960(%)\"\n$find.Replacement.Text = \"\ucd1d \uc774\uc775\ub960(%)\"\n$find.Replacement.Font.Bold = 1\n$find.Execute($find.Text, $false, $true, $false, $false, $false, $true, 1, $true, $find.Replacement.Text, 2)\n"}
